$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.131.04"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "1.788.39"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.85"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.92"
$ws.Range("E8").Value = "  -0.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.293"
$ws.Range("E9").Value = "  +1.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0688"
$ws.Range("E10").Value = "  -2.07%  "
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("D12").Value = "2.046.84"
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("E13").Value = "  +2.93%  "
$ws.Range("D14").Value = "1.790.05"
$ws.Range("E14").Value = "  +0.98%  "
$ws.Range("D15").Value = "34.069.74"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("E17").Value = "  +1.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.13"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.75"
$ws.Range("E19").Value = "  +1.34%  "
$ws.Range("D20").Value = "0.0₃0777"
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.87"
$ws.Range("E22").Value = "  +2.31%  "
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("E24").Value = "  -0.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.56"
$ws.Range("E25").Value = "  +0.98%  "
$ws.Range("E26").Value = "  +1.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.33"
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("E28").Value = "  +1.01%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  -0.68%  "
$ws.Range("E31").Value = "  +1.38%  "
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.61"
$ws.Range("E33").Value = "  +3.15%  "
$ws.Range("E34").Value = "  +0.35%  "
$ws.Range("D35").Value = "1.463.00"
$ws.Range("E35").Value = "  +5.20%  "
$ws.Range("E36").Value = "  +10.34%  "
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("E38").Value = "  +2.68%  "
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("E40").Value = "  +4.14%  "
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.920"
$ws.Range("E42").Value = "  +1.44%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.48"
$ws.Range("E44").Value = "  +3.08%  "
$ws.Range("E45").Value = "  +4.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0507"
$ws.Range("E46").Value = "  +2.03%  "
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("E48").Value = "  +0.73%  "
$ws.Range("D49").Value = "1.948.34"
$ws.Range("E49").Value = "  +0.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "106.34"
$ws.Range("E50").Value = "  -1.37%  "
$ws.Range("E51").Value = "  -0.07%  "
